$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2 is a zero-padded code stored as text ("002" -> "001"). Force text
# formatting first so Excel doesn't coerce it to the number 1, then drop
# the number format back to the sheet's default style.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "001"
$ws.Range("J2").Style = "Normal"

# N2 is a date-like value stored as text, not a real date serial.
$ws.Range("N2").NumberFormat = "@"
$ws.Range("N2").Value = "2018-12-31 00:00:00"
$ws.Range("N2").Style = "Normal"

# Numeric metrics for the new report period.
$ws.Range("O2").Value = 12763520.39
$ws.Range("P2").Value = 30.1421706372
$ws.Range("Q2").Value = 704845378.6
$ws.Range("R2").Value = 1664.5540591814
$ws.Range("S2").Value = 617464715.4400001
$ws.Range("T2").Value = 1458.1969743895
$ws.Range("U2").Value = 9311793.140000001
$ws.Range("V2").Value = 21.9906145944
$ws.Range("W2").Value = 1380506.48
$ws.Range("X2").Value = 3.2601868932
$ws.Range("Y2").Value = 6588381.84
$ws.Range("Z2").Value = 15.5590404196
$ws.Range("AA2").Value = 20269083.44
$ws.Range("AB2").Value = 47.8672147684
$ws.Range("AC2").Value = 42344396.97
$ws.Range("AD2").Value = 199.9085098946
